# Update the "Metric" text for the Molnupiravir / Paxlovid aged-care
# prescription rows to reflect that the figures are reported weekly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metrics")

$ws.Range("C72").Value = "# Aged Care Molnupiravir Prescriptions (Weekly)"
$ws.Range("C73").Value = "# Aged Care Molnupiravir Prescriptions (Weekly) per 1M"
$ws.Range("C74").Value = "% Aged Care Molnupiravir Prescriptions (Weekly) per Case"
$ws.Range("C75").Value = "# Aged Care Paxlovid Prescriptions (Weekly)"
$ws.Range("C76").Value = "# Aged Care Paxlovid Prescriptions (Weekly) per 1M"
$ws.Range("C77").Value = "% Aged Care Paxlovid Prescriptions (Weekly) per Case"

# Move the active selection to where the user left off after the edit.
$ws.Range("C78").Select()
